$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F18").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = 'Data Scientist / Data Science Specialist'
$ws.Cells.Item(2, 2).Value = 'Adidev Technologies Inc'
$ws.Cells.Item(2, 3).Value = 'Dallas, TX, US USA'
$ws.Cells.Item(2, 4).Value = 18.9
$ws.Cells.Item(2, 5).Value = 'Data Scientist, TensorFlow, PyTorch, AWS SageMaker, S3, EC2, MLflow, Jenkins, Git, Databricks'
$ws.Cells.Item(2, 6).Value = '2026-02-21'
$ws.Cells.Item(2, 7).Value = 'https://www.indeed.com/viewjob?jk=91cd87d39ed00187'

$ws.Cells.Item(3, 1).Value = 'Junior Software Engineer'
$ws.Cells.Item(3, 2).Value = 'MetLife'
$ws.Cells.Item(3, 3).Value = 'Whippany, NJ, US USA'
$ws.Cells.Item(3, 4).Value = 17.8
$ws.Cells.Item(3, 5).Value = 'LangChain, RAG, Copilot, Hugging Face, Prompt Engineering, TensorFlow, Docker, Kubernetes, CI/CD, Git'
$ws.Cells.Item(3, 6).Value = '2026-02-20'
$ws.Cells.Item(3, 7).Value = 'https://www.indeed.com/viewjob?jk=27a8a2a71007930a'

$ws.Cells.Item(4, 1).Value = 'Senior Machine Learning Engineer, Risk Modeling'
$ws.Cells.Item(4, 2).Value = 'Block'
$ws.Cells.Item(4, 3).Value = 'New York, NY, US USA'
$ws.Cells.Item(4, 4).Value = 15.6
$ws.Cells.Item(4, 5).Value = 'Machine Learning Engineer, TensorFlow, PyTorch, XGBoost, Keras, MLflow, CI/CD, Snowflake, PySpark, MySQL'
$ws.Cells.Item(4, 6).Value = '2026-02-20'
$ws.Cells.Item(4, 7).Value = 'https://www.indeed.com/viewjob?jk=2749ecb9f8a0dc09'

$ws.Cells.Item(5, 1).Value = 'Senior Machine Learning Engineer, Risk Modeling'
$ws.Cells.Item(5, 2).Value = 'Block'
$ws.Cells.Item(5, 3).Value = 'Los Angeles, CA, US USA'
$ws.Cells.Item(5, 4).Value = 15.6
$ws.Cells.Item(5, 5).Value = 'Machine Learning Engineer, TensorFlow, PyTorch, XGBoost, Keras, MLflow, CI/CD, Snowflake, PySpark, MySQL'
$ws.Cells.Item(5, 6).Value = '2026-02-20'
$ws.Cells.Item(5, 7).Value = 'https://www.indeed.com/viewjob?jk=43f57fcb65df08a3'

$ws.Cells.Item(6, 1).Value = 'Data Science'
$ws.Cells.Item(6, 2).Value = 'Adidev Technologies Inc'
$ws.Cells.Item(6, 3).Value = 'San Francisco, CA, US USA'
$ws.Cells.Item(6, 4).Value = 15.6
$ws.Cells.Item(6, 5).Value = 'Data Scientist, TensorFlow, PyTorch, XGBoost, Keras, spaCy, Kubernetes, PostgreSQL, MongoDB, Tableau'
$ws.Cells.Item(6, 6).Value = '2026-02-21'
$ws.Cells.Item(6, 7).Value = 'https://www.indeed.com/viewjob?jk=6fb5d815a109de6d'

$ws.Cells.Item(7, 1).Value = 'Data Science Specialist'
$ws.Cells.Item(7, 2).Value = 'Adidev Technologies Inc'
$ws.Cells.Item(7, 3).Value = 'New York, NY, US USA'
$ws.Cells.Item(7, 4).Value = 15.6
$ws.Cells.Item(7, 5).Value = 'Data Scientist, TensorFlow, PyTorch, XGBoost, Keras, spaCy, Kubernetes, PostgreSQL, MongoDB, Tableau'
$ws.Cells.Item(7, 6).Value = '2026-02-21'
$ws.Cells.Item(7, 7).Value = 'https://www.indeed.com/viewjob?jk=7686f3290697986d'

$ws.Cells.Item(8, 1).Value = 'DATA SCIENTIST'
$ws.Cells.Item(8, 2).Value = 'Adidev Technologies Inc'
$ws.Cells.Item(8, 3).Value = 'Austin, TX, US USA'
$ws.Cells.Item(8, 4).Value = 15.6
$ws.Cells.Item(8, 5).Value = 'Data Scientist, TensorFlow, PyTorch, XGBoost, Keras, spaCy, Kubernetes, PostgreSQL, MongoDB, Tableau'
$ws.Cells.Item(8, 6).Value = '2026-02-21'
$ws.Cells.Item(8, 7).Value = 'https://www.indeed.com/viewjob?jk=2b438984f6dd9115'

$ws.Cells.Item(9, 1).Value = 'Software Engineer - Full Stack'
$ws.Cells.Item(9, 2).Value = 'NTT DATA'
$ws.Cells.Item(9, 3).Value = 'MO, US USA'
$ws.Cells.Item(9, 4).Value = 14.4
$ws.Cells.Item(9, 5).Value = 'RAG, Docker, Kubernetes, CI/CD, Jenkins, Git, PostgreSQL, MongoDB, NoSQL, SQL'
$ws.Cells.Item(9, 6).Value = '2026-02-20'
$ws.Cells.Item(9, 7).Value = 'https://www.indeed.com/viewjob?jk=e06ca5d0cbf6168e'

$ws.Cells.Item(10, 1).Value = 'Dev Ops and Cloud Engineer, Associate'
$ws.Cells.Item(10, 2).Value = 'BlackRock'
$ws.Cells.Item(10, 3).Value = 'Atlanta, GA, US USA'
$ws.Cells.Item(10, 4).Value = 13.3
$ws.Cells.Item(10, 5).Value = 'Copilot, Docker, Kubernetes, AKS, CI/CD, GitHub Actions, Terraform, Git, Python, R'
$ws.Cells.Item(10, 6).Value = '2026-02-20'
$ws.Cells.Item(10, 7).Value = 'https://www.indeed.com/viewjob?jk=e73c18f9abaa746c'

$ws.Cells.Item(11, 1).Value = 'Data Engineer III'
$ws.Cells.Item(11, 2).Value = 'Grainger'
$ws.Cells.Item(11, 3).Value = 'Chicago, IL, US USA'
$ws.Cells.Item(11, 4).Value = 13.3
$ws.Cells.Item(11, 5).Value = 'Data Scientist, RAG, Docker, Kubernetes, CI/CD, Git, Snowflake, Kafka, Python, SQL'
$ws.Cells.Item(11, 6).Value = '2026-02-20'
$ws.Cells.Item(11, 7).Value = 'https://www.indeed.com/viewjob?jk=ca71348ca74ae1ce'

$ws.Cells.Item(12, 1).Value = 'Software Engineer'
$ws.Cells.Item(12, 2).Value = 'McKesson'
$ws.Cells.Item(12, 3).Value = 'Columbus, OH, US USA'
$ws.Cells.Item(12, 4).Value = 12.2
$ws.Cells.Item(12, 5).Value = 'RAG, CI/CD, Terraform, Git, Kafka, PostgreSQL, SQL, R, Java, Scala'
$ws.Cells.Item(12, 6).Value = '2026-02-20'
$ws.Cells.Item(12, 7).Value = 'https://www.indeed.com/viewjob?jk=a1ace665b571d0bb'

$ws.Cells.Item(13, 1).Value = 'Sr Software Engineer'
$ws.Cells.Item(13, 2).Value = 'McKesson'
$ws.Cells.Item(13, 3).Value = 'Columbus, OH, US USA'
$ws.Cells.Item(13, 4).Value = 12.2
$ws.Cells.Item(13, 5).Value = 'RAG, CI/CD, Terraform, Git, Kafka, PostgreSQL, SQL, R, Java, Scala'
$ws.Cells.Item(13, 6).Value = '2026-02-20'
$ws.Cells.Item(13, 7).Value = 'https://www.indeed.com/viewjob?jk=75d8c3155721fa1b'

$ws.Cells.Item(14, 1).Value = 'Data Platform Engineer II'
$ws.Cells.Item(14, 2).Value = 'Best Egg'
$ws.Cells.Item(14, 3).Value = 'Wilmington, DE, US USA'
$ws.Cells.Item(14, 4).Value = 12.2
$ws.Cells.Item(14, 5).Value = 'RAG, S3, Data Lake, Docker, CI/CD, Git, Snowflake, Python, SQL, R'
$ws.Cells.Item(14, 6).Value = '2026-02-20'
$ws.Cells.Item(14, 7).Value = 'https://www.indeed.com/viewjob?jk=1ff6ad916136713c'

$ws.Cells.Item(15, 1).Value = 'Sr. Quality Engineer'
$ws.Cells.Item(15, 2).Value = 'McKesson'
$ws.Cells.Item(15, 3).Value = 'Columbus, OH, US USA'
$ws.Cells.Item(15, 4).Value = 11.1
$ws.Cells.Item(15, 5).Value = 'RAG, CI/CD, GitHub Actions, Git, Kafka, PostgreSQL, SQL, R, Java, Scala'
$ws.Cells.Item(15, 6).Value = '2026-02-20'
$ws.Cells.Item(15, 7).Value = 'https://www.indeed.com/viewjob?jk=484a7d4e140ae699'

$ws.Cells.Item(16, 1).Value = 'Senior Software Engineer - AI Research Clusters'
$ws.Cells.Item(16, 2).Value = 'NVIDIA'
$ws.Cells.Item(16, 3).Value = 'Santa Clara, CA, US USA'
$ws.Cells.Item(16, 4).Value = 11.1
$ws.Cells.Item(16, 5).Value = 'Generative AI, RAG, Cortex, Docker, Kubernetes, Git, Python, R, Java, Optimization'
$ws.Cells.Item(16, 6).Value = '2026-02-20'
$ws.Cells.Item(16, 7).Value = 'https://www.indeed.com/viewjob?jk=e62aae990e8e1e0d'

$ws.Cells.Item(17, 1).Value = 'AI/ML Engineer SME (CMS)'
$ws.Cells.Item(17, 2).Value = 'General Dynamics Information Technology'
$ws.Cells.Item(17, 3).Value = 'IN, US USA'
$ws.Cells.Item(17, 4).Value = 11.1
$ws.Cells.Item(17, 5).Value = 'Data Scientist, Generative AI, Git, Snowflake, Databricks, Tableau, Python, SQL, R, Optimization'
$ws.Cells.Item(17, 6).Value = '2026-02-20'
$ws.Cells.Item(17, 7).Value = 'https://www.indeed.com/viewjob?jk=fdbbee24455a9dd6'

$ws.Cells.Item(18, 1).Value = 'AI and ML HPC Cluster Engineer'
$ws.Cells.Item(18, 2).Value = 'NVIDIA'
$ws.Cells.Item(18, 3).Value = 'Santa Clara, CA, US USA'
$ws.Cells.Item(18, 4).Value = 10
$ws.Cells.Item(18, 5).Value = 'Generative AI, RAG, TensorFlow, PyTorch, Docker, Python, R, Scala, Optimization'
$ws.Cells.Item(18, 6).Value = '2026-02-20'
$ws.Cells.Item(18, 7).Value = 'https://www.indeed.com/viewjob?jk=608e4534825362a1'

$ws.Range("F2:F18").Style = "Normal"
